$wb = $excel.ActiveWorkbook

# LL_max_6 (sheet9.xml)
$ws = $wb.Worksheets.Item("LL_max_6")
$ws.Range("O3").Value = [double]"0.291705603751294"
$ws.Range("Q3").Value = [double]"7.091324389062566"
$ws.Range("C4").Value = [double]"3.709154423937071"
$ws.Range("D4").Value = [double]"3.709154423937071"
$ws.Range("F4").Value = [double]"42.82962610251919"
$ws.Range("G4").Value = [double]"42.82962610251919"
$ws.Range("H4").Value = [double]"4.37704938291617"
$ws.Range("I4").Value = [double]"3.22760593527366"
$ws.Range("J4").Value = [double]"0.6277319163569961"
$ws.Range("K4").Value = [double]"2.898438697948067"
$ws.Range("L4").Value = [double]"0.627731916364051"
$ws.Range("M4").Value = [double]"2.898438697960011"
$ws.Range("N4").Value = [double]"0.9526279648032089"
$ws.Range("P4").Value = [double]"0.9526279648047555"
$ws.Range("Q4").Value = [double]"-1.2483317126697E-11"
$ws.Range("N5").Value = [double]"0.9526279647942851"
$ws.Range("Q5").Value = [double]"4.605544102973149E-10"
$ws.Range("P6").Value = [double]"0.9526279648136791"
$ws.Range("Q6").Value = [double]"4.605544102973149E-10"

# LL_max_10 (sheet10.xml)
$ws = $wb.Worksheets.Item("LL_max_10")
$ws.Range("O3").Value = [double]"0.291705603751294"
$ws.Range("Q3").Value = [double]"7.091324389062566"
$ws.Range("C4").Value = [double]"3.709154423937071"
$ws.Range("D4").Value = [double]"3.709154423937071"
$ws.Range("F4").Value = [double]"42.82962610251919"
$ws.Range("G4").Value = [double]"42.82962610251919"
$ws.Range("H4").Value = [double]"4.37704938291617"
$ws.Range("I4").Value = [double]"3.22760593527366"
$ws.Range("J4").Value = [double]"0.6277319163569961"
$ws.Range("K4").Value = [double]"2.898438697948067"
$ws.Range("L4").Value = [double]"0.627731916364051"
$ws.Range("M4").Value = [double]"2.898438697960011"
$ws.Range("N4").Value = [double]"0.9526279648032089"
$ws.Range("P4").Value = [double]"0.9526279648047555"
$ws.Range("Q4").Value = [double]"-1.2483317126697E-11"
$ws.Range("N5").Value = [double]"0.9526279647942851"
$ws.Range("Q5").Value = [double]"4.605544102973149E-10"
$ws.Range("P6").Value = [double]"0.9526279648136791"
$ws.Range("Q6").Value = [double]"4.605544102973149E-10"

# LL_max_fault_6 (sheet11.xml)
$ws = $wb.Worksheets.Item("LL_max_fault_6")
$ws.Range("O3").Value = [double]"0.291705603751294"
$ws.Range("Q3").Value = [double]"7.091324389062566"
$ws.Range("C4").Value = [double]"3.709154423937071"
$ws.Range("D4").Value = [double]"3.709154423937071"
$ws.Range("F4").Value = [double]"42.82962610251919"
$ws.Range("G4").Value = [double]"42.82962610251919"
$ws.Range("H4").Value = [double]"4.37704938291617"
$ws.Range("I4").Value = [double]"3.22760593527366"
$ws.Range("J4").Value = [double]"0.6277319163569961"
$ws.Range("K4").Value = [double]"2.898438697948067"
$ws.Range("L4").Value = [double]"0.627731916364051"
$ws.Range("M4").Value = [double]"2.898438697960011"
$ws.Range("N4").Value = [double]"0.9526279648032089"
$ws.Range("P4").Value = [double]"0.9526279648047555"
$ws.Range("Q4").Value = [double]"-1.2483317126697E-11"
$ws.Range("N5").Value = [double]"0.9526279647942851"
$ws.Range("Q5").Value = [double]"4.605544102973149E-10"
$ws.Range("P6").Value = [double]"0.9526279648136791"
$ws.Range("Q6").Value = [double]"4.605544102973149E-10"

# LL_max_fault_10 (sheet12.xml)
$ws = $wb.Worksheets.Item("LL_max_fault_10")
$ws.Range("O3").Value = [double]"0.291705603751294"
$ws.Range("Q3").Value = [double]"7.091324389062566"
$ws.Range("C4").Value = [double]"3.709154423937071"
$ws.Range("D4").Value = [double]"3.709154423937071"
$ws.Range("F4").Value = [double]"42.82962610251919"
$ws.Range("G4").Value = [double]"42.82962610251919"
$ws.Range("H4").Value = [double]"4.37704938291617"
$ws.Range("I4").Value = [double]"3.22760593527366"
$ws.Range("J4").Value = [double]"0.6277319163569961"
$ws.Range("K4").Value = [double]"2.898438697948067"
$ws.Range("L4").Value = [double]"0.627731916364051"
$ws.Range("M4").Value = [double]"2.898438697960011"
$ws.Range("N4").Value = [double]"0.9526279648032089"
$ws.Range("P4").Value = [double]"0.9526279648047555"
$ws.Range("Q4").Value = [double]"-1.2483317126697E-11"
$ws.Range("N5").Value = [double]"0.9526279647942851"
$ws.Range("Q5").Value = [double]"4.605544102973149E-10"
$ws.Range("P6").Value = [double]"0.9526279648136791"
$ws.Range("Q6").Value = [double]"4.605544102973149E-10"

# LL_min_6 (sheet13.xml)
$ws = $wb.Worksheets.Item("LL_min_6")
$ws.Range("O3").Value = [double]"0.3617201835030113"
$ws.Range("P3").Value = [double]"0.7649715466411117"
$ws.Range("H4").Value = [double]"6.086106992372931"
$ws.Range("I4").Value = [double]"3.212307515692304"
$ws.Range("J4").Value = [double]"1.113751880246063"
$ws.Range("K4").Value = [double]"3.006397543767953"
$ws.Range("L4").Value = [double]"1.113751880233333"
$ws.Range("M4").Value = [double]"3.006397543745581"
$ws.Range("P4").Value = [double]"0.8660254037883345"
$ws.Range("Q4").Value = [double]"1.723409298412112E-10"
$ws.Range("N5").Value = [double]"0.8660254037798931"
$ws.Range("P5").Value = [double]"0.8660254037964471"
$ws.Range("Q5").Value = [double]"1.080573051089087E-09"
$ws.Range("N6").Value = [double]"0.8660254037798931"
$ws.Range("P6").Value = [double]"0.8660254037964471"
$ws.Range("Q6").Value = [double]"1.080573051089087E-09"

# LL_min_10 (sheet14.xml)
$ws = $wb.Worksheets.Item("LL_min_10")
$ws.Range("O3").Value = [double]"0.3617201835030113"
$ws.Range("P3").Value = [double]"0.7649715466411117"
$ws.Range("H4").Value = [double]"6.086106992372931"
$ws.Range("I4").Value = [double]"3.212307515692304"
$ws.Range("J4").Value = [double]"1.113751880246063"
$ws.Range("K4").Value = [double]"3.006397543767953"
$ws.Range("L4").Value = [double]"1.113751880233333"
$ws.Range("M4").Value = [double]"3.006397543745581"
$ws.Range("P4").Value = [double]"0.8660254037883345"
$ws.Range("Q4").Value = [double]"1.723409298412112E-10"
$ws.Range("N5").Value = [double]"0.8660254037798931"
$ws.Range("P5").Value = [double]"0.8660254037964471"
$ws.Range("Q5").Value = [double]"1.080573051089087E-09"
$ws.Range("N6").Value = [double]"0.8660254037798931"
$ws.Range("P6").Value = [double]"0.8660254037964471"
$ws.Range("Q6").Value = [double]"1.080573051089087E-09"

# LL_min_fault_6 (sheet15.xml)
$ws = $wb.Worksheets.Item("LL_min_fault_6")
$ws.Range("O3").Value = [double]"0.3617201835030113"
$ws.Range("P3").Value = [double]"0.7649715466411117"
$ws.Range("H4").Value = [double]"6.086106992372931"
$ws.Range("I4").Value = [double]"3.212307515692304"
$ws.Range("J4").Value = [double]"1.113751880246063"
$ws.Range("K4").Value = [double]"3.006397543767953"
$ws.Range("L4").Value = [double]"1.113751880233333"
$ws.Range("M4").Value = [double]"3.006397543745581"
$ws.Range("P4").Value = [double]"0.8660254037883345"
$ws.Range("Q4").Value = [double]"1.723409298412112E-10"
$ws.Range("N5").Value = [double]"0.8660254037798931"
$ws.Range("P5").Value = [double]"0.8660254037964471"
$ws.Range("Q5").Value = [double]"1.080573051089087E-09"
$ws.Range("N6").Value = [double]"0.8660254037798931"
$ws.Range("P6").Value = [double]"0.8660254037964471"
$ws.Range("Q6").Value = [double]"1.080573051089087E-09"

# LL_min_fault_10 (sheet16.xml)
$ws = $wb.Worksheets.Item("LL_min_fault_10")
$ws.Range("O3").Value = [double]"0.3617201835030113"
$ws.Range("P3").Value = [double]"0.7649715466411117"
$ws.Range("H4").Value = [double]"6.086106992372931"
$ws.Range("I4").Value = [double]"3.212307515692304"
$ws.Range("J4").Value = [double]"1.113751880246063"
$ws.Range("K4").Value = [double]"3.006397543767953"
$ws.Range("L4").Value = [double]"1.113751880233333"
$ws.Range("M4").Value = [double]"3.006397543745581"
$ws.Range("P4").Value = [double]"0.8660254037883345"
$ws.Range("Q4").Value = [double]"1.723409298412112E-10"
$ws.Range("N5").Value = [double]"0.8660254037798931"
$ws.Range("P5").Value = [double]"0.8660254037964471"
$ws.Range("Q5").Value = [double]"1.080573051089087E-09"
$ws.Range("N6").Value = [double]"0.8660254037798931"
$ws.Range("P6").Value = [double]"0.8660254037964471"
$ws.Range("Q6").Value = [double]"1.080573051089087E-09"

# LG_max_6 (sheet17.xml)
$ws = $wb.Worksheets.Item("LG_max_6")
$ws.Range("N2").Value = [double]"1.021377084239933"
$ws.Range("O2").Value = [double]"1.10000002384589"
$ws.Range("P2").Value = [double]"1.058260365617711"
$ws.Range("Q2").Value = [double]"30.28841431342822"
$ws.Range("R2").Value = [double]"-89.99999999999633"
$ws.Range("S2").Value = [double]"146.449838450249"
$ws.Range("N3").Value = [double]"0.68615965300225"
$ws.Range("O3").Value = [double]"1.100000023849745"
$ws.Range("P3").Value = [double]"0.9572858912279797"
$ws.Range("Q3").Value = [double]"30.42354294779994"
$ws.Range("R3").Value = [double]"-89.99999999999706"
$ws.Range("S3").Value = [double]"128.1760681478563"
$ws.Range("B4").Value = [double]"3.581989981078915"
$ws.Range("E4").Value = [double]"41.36125759620908"
$ws.Range("H4").Value = [double]"4.377049382916195"
$ws.Range("I4").Value = [double]"3.227605935273422"
$ws.Range("J4").Value = [double]"0.6277319162689007"
$ws.Range("K4").Value = [double]"2.898438697960076"
$ws.Range("L4").Value = [double]"0.6277319163639457"
$ws.Range("M4").Value = [double]"2.898438697959814"
$ws.Range("N4").Value = [double]"0.5850190725743573"
$ws.Range("O4").Value = [double]"1.100000023849745"
$ws.Range("P4").Value = [double]"0.8585464144322896"
$ws.Range("Q4").Value = [double]"39.29756343703533"
$ws.Range("R4").Value = [double]"-89.99999999999714"
$ws.Range("S4").Value = [double]"121.8244261349591"
$ws.Range("T4").Value = [double]"3.581989981078914"
$ws.Range("N5").Value = [double]"0.5850190725844024"
$ws.Range("O5").Value = [double]"1.100000023849745"
$ws.Range("P5").Value = [double]"0.8585464144254525"
$ws.Range("Q5").Value = [double]"39.29756343758146"
$ws.Range("R5").Value = [double]"-89.99999999999727"
$ws.Range("S5").Value = [double]"121.8244261355752"
$ws.Range("N6").Value = [double]"0.5850190725844024"
$ws.Range("O6").Value = [double]"1.100000023849745"
$ws.Range("P6").Value = [double]"0.8585464144254525"
$ws.Range("Q6").Value = [double]"39.29756343758146"
$ws.Range("R6").Value = [double]"-89.99999999999727"
$ws.Range("S6").Value = [double]"121.8244261355752"

# LG_max_10 (sheet18.xml)
$ws = $wb.Worksheets.Item("LG_max_10")
$ws.Range("N2").Value = [double]"1.021377084239933"
$ws.Range("O2").Value = [double]"1.10000002384589"
$ws.Range("P2").Value = [double]"1.058260365617711"
$ws.Range("Q2").Value = [double]"30.28841431342822"
$ws.Range("R2").Value = [double]"-89.99999999999633"
$ws.Range("S2").Value = [double]"146.449838450249"
$ws.Range("N3").Value = [double]"0.68615965300225"
$ws.Range("O3").Value = [double]"1.100000023849745"
$ws.Range("P3").Value = [double]"0.9572858912279797"
$ws.Range("Q3").Value = [double]"30.42354294779994"
$ws.Range("R3").Value = [double]"-89.99999999999706"
$ws.Range("S3").Value = [double]"128.1760681478563"
$ws.Range("B4").Value = [double]"3.581989981078915"
$ws.Range("E4").Value = [double]"41.36125759620908"
$ws.Range("H4").Value = [double]"4.377049382916195"
$ws.Range("I4").Value = [double]"3.227605935273422"
$ws.Range("J4").Value = [double]"0.6277319162689007"
$ws.Range("K4").Value = [double]"2.898438697960076"
$ws.Range("L4").Value = [double]"0.6277319163639457"
$ws.Range("M4").Value = [double]"2.898438697959814"
$ws.Range("N4").Value = [double]"0.5850190725743573"
$ws.Range("O4").Value = [double]"1.100000023849745"
$ws.Range("P4").Value = [double]"0.8585464144322896"
$ws.Range("Q4").Value = [double]"39.29756343703533"
$ws.Range("R4").Value = [double]"-89.99999999999714"
$ws.Range("S4").Value = [double]"121.8244261349591"
$ws.Range("T4").Value = [double]"3.581989981078914"
$ws.Range("N5").Value = [double]"0.5850190725844024"
$ws.Range("O5").Value = [double]"1.100000023849745"
$ws.Range("P5").Value = [double]"0.8585464144254525"
$ws.Range("Q5").Value = [double]"39.29756343758146"
$ws.Range("R5").Value = [double]"-89.99999999999727"
$ws.Range("S5").Value = [double]"121.8244261355752"
$ws.Range("N6").Value = [double]"0.5850190725844024"
$ws.Range("O6").Value = [double]"1.100000023849745"
$ws.Range("P6").Value = [double]"0.8585464144254525"
$ws.Range("Q6").Value = [double]"39.29756343758146"
$ws.Range("R6").Value = [double]"-89.99999999999727"
$ws.Range("S6").Value = [double]"121.8244261355752"

# LG_max_fault_6 (sheet19.xml)
$ws = $wb.Worksheets.Item("LG_max_fault_6")
$ws.Range("N2").Value = [double]"1.073629283973476"
$ws.Range("O2").Value = [double]"1.100000023843262"
$ws.Range("P2").Value = [double]"1.089053953701369"
$ws.Range("Q2").Value = [double]"29.87820653143294"
$ws.Range("S2").Value = [double]"148.7382466549213"
$ws.Range("N3").Value = [double]"0.9625077933728241"
$ws.Range("O3").Value = [double]"1.100000023844815"
$ws.Range("P3").Value = [double]"1.055779646083895"
$ws.Range("Q3").Value = [double]"28.85030996543065"
$ws.Range("R3").Value = [double]"-89.99999999999666"
$ws.Range("S3").Value = [double]"142.9875046145783"
$ws.Range("B4").Value = [double]"1.203260365654794"
$ws.Range("E4").Value = [double]"13.89405392032006"
$ws.Range("H4").Value = [double]"4.377049382916195"
$ws.Range("I4").Value = [double]"3.227605935273422"
$ws.Range("J4").Value = [double]"0.6277319162689007"
$ws.Range("K4").Value = [double]"2.898438697960076"
$ws.Range("L4").Value = [double]"0.6277319163639457"
$ws.Range("M4").Value = [double]"2.898438697959814"
$ws.Range("N4").Value = [double]"0.921657749219949"
$ws.Range("O4").Value = [double]"1.100000023844815"
$ws.Range("P4").Value = [double]"1.01802104518981"
$ws.Range("Q4").Value = [double]"30.30269163714974"
$ws.Range("R4").Value = [double]"-89.99999999999667"
$ws.Range("S4").Value = [double]"141.4116448044574"
$ws.Range("T4").Value = [double]"1.203260365654794"
$ws.Range("N5").Value = [double]"0.9216577492225949"
$ws.Range("O5").Value = [double]"1.100000023844815"
$ws.Range("P5").Value = [double]"1.018021045188142"
$ws.Range("Q5").Value = [double]"30.30269163732438"
$ws.Range("R5").Value = [double]"-89.99999999999672"
$ws.Range("S5").Value = [double]"141.4116448046532"
$ws.Range("N6").Value = [double]"0.9216577492225949"
$ws.Range("O6").Value = [double]"1.100000023844815"
$ws.Range("P6").Value = [double]"1.018021045188142"
$ws.Range("Q6").Value = [double]"30.30269163732438"
$ws.Range("R6").Value = [double]"-89.99999999999672"
$ws.Range("S6").Value = [double]"141.4116448046532"

# LG_max_fault_10 (sheet20.xml)
$ws = $wb.Worksheets.Item("LG_max_fault_10")
$ws.Range("N2").Value = [double]"1.073629283973476"
$ws.Range("O2").Value = [double]"1.100000023843262"
$ws.Range("P2").Value = [double]"1.089053953701369"
$ws.Range("Q2").Value = [double]"29.87820653143294"
$ws.Range("S2").Value = [double]"148.7382466549213"
$ws.Range("N3").Value = [double]"0.9625077933728241"
$ws.Range("O3").Value = [double]"1.100000023844815"
$ws.Range("P3").Value = [double]"1.055779646083895"
$ws.Range("Q3").Value = [double]"28.85030996543065"
$ws.Range("R3").Value = [double]"-89.99999999999666"
$ws.Range("S3").Value = [double]"142.9875046145783"
$ws.Range("B4").Value = [double]"1.203260365654794"
$ws.Range("E4").Value = [double]"13.89405392032006"
$ws.Range("H4").Value = [double]"4.377049382916195"
$ws.Range("I4").Value = [double]"3.227605935273422"
$ws.Range("J4").Value = [double]"0.6277319162689007"
$ws.Range("K4").Value = [double]"2.898438697960076"
$ws.Range("L4").Value = [double]"0.6277319163639457"
$ws.Range("M4").Value = [double]"2.898438697959814"
$ws.Range("N4").Value = [double]"0.921657749219949"
$ws.Range("O4").Value = [double]"1.100000023844815"
$ws.Range("P4").Value = [double]"1.01802104518981"
$ws.Range("Q4").Value = [double]"30.30269163714974"
$ws.Range("R4").Value = [double]"-89.99999999999667"
$ws.Range("S4").Value = [double]"141.4116448044574"
$ws.Range("T4").Value = [double]"1.203260365654794"
$ws.Range("N5").Value = [double]"0.9216577492225949"
$ws.Range("O5").Value = [double]"1.100000023844815"
$ws.Range("P5").Value = [double]"1.018021045188142"
$ws.Range("Q5").Value = [double]"30.30269163732438"
$ws.Range("R5").Value = [double]"-89.99999999999672"
$ws.Range("S5").Value = [double]"141.4116448046532"
$ws.Range("N6").Value = [double]"0.9216577492225949"
$ws.Range("O6").Value = [double]"1.100000023844815"
$ws.Range("P6").Value = [double]"1.018021045188142"
$ws.Range("Q6").Value = [double]"30.30269163732438"
$ws.Range("R6").Value = [double]"-89.99999999999672"
$ws.Range("S6").Value = [double]"141.4116448046532"

# LG_min_6 (sheet21.xml)
$ws = $wb.Worksheets.Item("LG_min_6")
$ws.Range("N2").Value = [double]"0.9307178477919329"
$ws.Range("O2").Value = [double]"0.9999999999965741"
$ws.Range("P2").Value = [double]"0.9740812140644213"
$ws.Range("Q2").Value = [double]"29.52782603558825"
$ws.Range("R2").Value = [double]"-89.99999999999636"
$ws.Range("S2").Value = [double]"146.2408933423334"
$ws.Range("N3").Value = [double]"0.6692408528128602"
$ws.Range("O3").Value = [double]"1.000000000000118"
$ws.Range("P3").Value = [double]"0.9355345379569088"
$ws.Range("Q3").Value = [double]"25.33067431500087"
$ws.Range("R3").Value = [double]"-89.99999999999724"
$ws.Range("S3").Value = [double]"130.2840694740472"
$ws.Range("B4").Value = [double]"2.789473153659462"
$ws.Range("E4").Value = [double]"32.21006152325049"
$ws.Range("H4").Value = [double]"6.086106992372986"
$ws.Range("I4").Value = [double]"3.212307515692044"
$ws.Range("J4").Value = [double]"1.113751880275147"
$ws.Range("K4").Value = [double]"3.006397543745625"
$ws.Range("L4").Value = [double]"1.113751880233235"
$ws.Range("M4").Value = [double]"3.006397543745307"
$ws.Range("N4").Value = [double]"0.5612515215730326"
$ws.Range("O4").Value = [double]"1.000000000000118"
$ws.Range("P4").Value = [double]"0.8035158449962232"
$ws.Range("Q4").Value = [double]"36.60644220974014"
$ws.Range("R4").Value = [double]"-89.99999999999731"
$ws.Range("S4").Value = [double]"124.1053896004218"
$ws.Range("T4").Value = [double]"2.789473153659462"
$ws.Range("N5").Value = [double]"0.5612515215867692"
$ws.Range("O5").Value = [double]"1.000000000000118"
$ws.Range("P5").Value = [double]"0.8035158449917778"
$ws.Range("Q5").Value = [double]"36.60644221013293"
$ws.Range("R5").Value = [double]"-89.99999999999753"
$ws.Range("S5").Value = [double]"124.1053896013882"
$ws.Range("N6").Value = [double]"0.5612515215867692"
$ws.Range("O6").Value = [double]"1.000000000000118"
$ws.Range("P6").Value = [double]"0.8035158449917778"
$ws.Range("Q6").Value = [double]"36.60644221013293"
$ws.Range("R6").Value = [double]"-89.99999999999753"
$ws.Range("S6").Value = [double]"124.1053896013882"

# LG_min_10 (sheet22.xml)
$ws = $wb.Worksheets.Item("LG_min_10")
$ws.Range("N2").Value = [double]"0.9307178477919329"
$ws.Range("O2").Value = [double]"0.9999999999965741"
$ws.Range("P2").Value = [double]"0.9740812140644213"
$ws.Range("Q2").Value = [double]"29.52782603558825"
$ws.Range("R2").Value = [double]"-89.99999999999636"
$ws.Range("S2").Value = [double]"146.2408933423334"
$ws.Range("N3").Value = [double]"0.6692408528128602"
$ws.Range("O3").Value = [double]"1.000000000000118"
$ws.Range("P3").Value = [double]"0.9355345379569088"
$ws.Range("Q3").Value = [double]"25.33067431500087"
$ws.Range("R3").Value = [double]"-89.99999999999724"
$ws.Range("S3").Value = [double]"130.2840694740472"
$ws.Range("B4").Value = [double]"2.789473153659462"
$ws.Range("E4").Value = [double]"32.21006152325049"
$ws.Range("H4").Value = [double]"6.086106992372986"
$ws.Range("I4").Value = [double]"3.212307515692044"
$ws.Range("J4").Value = [double]"1.113751880275147"
$ws.Range("K4").Value = [double]"3.006397543745625"
$ws.Range("L4").Value = [double]"1.113751880233235"
$ws.Range("M4").Value = [double]"3.006397543745307"
$ws.Range("N4").Value = [double]"0.5612515215730326"
$ws.Range("O4").Value = [double]"1.000000000000118"
$ws.Range("P4").Value = [double]"0.8035158449962232"
$ws.Range("Q4").Value = [double]"36.60644220974014"
$ws.Range("R4").Value = [double]"-89.99999999999731"
$ws.Range("S4").Value = [double]"124.1053896004218"
$ws.Range("T4").Value = [double]"2.789473153659462"
$ws.Range("N5").Value = [double]"0.5612515215867692"
$ws.Range("O5").Value = [double]"1.000000000000118"
$ws.Range("P5").Value = [double]"0.8035158449917778"
$ws.Range("Q5").Value = [double]"36.60644221013293"
$ws.Range("R5").Value = [double]"-89.99999999999753"
$ws.Range("S5").Value = [double]"124.1053896013882"
$ws.Range("N6").Value = [double]"0.5612515215867692"
$ws.Range("O6").Value = [double]"1.000000000000118"
$ws.Range("P6").Value = [double]"0.8035158449917778"
$ws.Range("Q6").Value = [double]"36.60644221013293"
$ws.Range("R6").Value = [double]"-89.99999999999753"
$ws.Range("S6").Value = [double]"124.1053896013882"

# LG_min_fault_6 (sheet23.xml)
$ws = $wb.Worksheets.Item("LG_min_fault_6")
$ws.Range("N2").Value = [double]"0.9745113091810109"
$ws.Range("O2").Value = [double]"0.999999999998688"
$ws.Range("P2").Value = [double]"0.9907190801093688"
$ws.Range("Q2").Value = [double]"29.78423173290928"
$ws.Range("S2").Value = [double]"148.6168140650805"
$ws.Range("N3").Value = [double]"0.8780654026836471"
$ws.Range("O3").Value = [double]"0.9999999999999872"
$ws.Range("P3").Value = [double]"0.9669934339145485"
$ws.Range("Q3").Value = [double]"28.42463894915672"
$ws.Range("R3").Value = [double]"-89.99999999999669"
$ws.Range("S3").Value = [double]"142.9935043780143"
$ws.Range("B4").Value = [double]"1.03033483292163"
$ws.Range("E4").Value = [double]"11.89728186285502"
$ws.Range("H4").Value = [double]"6.086106992372986"
$ws.Range("I4").Value = [double]"3.212307515692044"
$ws.Range("J4").Value = [double]"1.113751880275147"
$ws.Range("K4").Value = [double]"3.006397543745625"
$ws.Range("L4").Value = [double]"1.113751880233235"
$ws.Range("M4").Value = [double]"3.006397543745307"
$ws.Range("N4").Value = [double]"0.8360759682777967"
$ws.Range("O4").Value = [double]"0.9999999999999872"
$ws.Range("P4").Value = [double]"0.9116650907700807"
$ws.Range("Q4").Value = [double]"31.26691834742913"
$ws.Range("R4").Value = [double]"-89.99999999999672"
$ws.Range("S4").Value = [double]"141.6179324090518"
$ws.Range("T4").Value = [double]"1.03033483292163"
$ws.Range("N5").Value = [double]"0.8360759682826518"
$ws.Range("O5").Value = [double]"0.9999999999999872"
$ws.Range("P5").Value = [double]"0.9116650907698488"
$ws.Range("Q5").Value = [double]"31.26691834756943"
$ws.Range("R5").Value = [double]"-89.9999999999968"
$ws.Range("S5").Value = [double]"141.6179324093825"
$ws.Range("N6").Value = [double]"0.8360759682826518"
$ws.Range("O6").Value = [double]"0.9999999999999872"
$ws.Range("P6").Value = [double]"0.9116650907698488"
$ws.Range("Q6").Value = [double]"31.26691834756943"
$ws.Range("R6").Value = [double]"-89.9999999999968"
$ws.Range("S6").Value = [double]"141.6179324093825"

# LG_min_fault_10 (sheet24.xml)
$ws = $wb.Worksheets.Item("LG_min_fault_10")
$ws.Range("N2").Value = [double]"0.9745113091810109"
$ws.Range("O2").Value = [double]"0.999999999998688"
$ws.Range("P2").Value = [double]"0.9907190801093688"
$ws.Range("Q2").Value = [double]"29.78423173290928"
$ws.Range("S2").Value = [double]"148.6168140650805"
$ws.Range("N3").Value = [double]"0.8780654026836471"
$ws.Range("O3").Value = [double]"0.9999999999999872"
$ws.Range("P3").Value = [double]"0.9669934339145485"
$ws.Range("Q3").Value = [double]"28.42463894915672"
$ws.Range("R3").Value = [double]"-89.99999999999669"
$ws.Range("S3").Value = [double]"142.9935043780143"
$ws.Range("B4").Value = [double]"1.03033483292163"
$ws.Range("E4").Value = [double]"11.89728186285502"
$ws.Range("H4").Value = [double]"6.086106992372986"
$ws.Range("I4").Value = [double]"3.212307515692044"
$ws.Range("J4").Value = [double]"1.113751880275147"
$ws.Range("K4").Value = [double]"3.006397543745625"
$ws.Range("L4").Value = [double]"1.113751880233235"
$ws.Range("M4").Value = [double]"3.006397543745307"
$ws.Range("N4").Value = [double]"0.8360759682777967"
$ws.Range("O4").Value = [double]"0.9999999999999872"
$ws.Range("P4").Value = [double]"0.9116650907700807"
$ws.Range("Q4").Value = [double]"31.26691834742913"
$ws.Range("R4").Value = [double]"-89.99999999999672"
$ws.Range("S4").Value = [double]"141.6179324090518"
$ws.Range("T4").Value = [double]"1.03033483292163"
$ws.Range("N5").Value = [double]"0.8360759682826518"
$ws.Range("O5").Value = [double]"0.9999999999999872"
$ws.Range("P5").Value = [double]"0.9116650907698488"
$ws.Range("Q5").Value = [double]"31.26691834756943"
$ws.Range("R5").Value = [double]"-89.9999999999968"
$ws.Range("S5").Value = [double]"141.6179324093825"
$ws.Range("N6").Value = [double]"0.8360759682826518"
$ws.Range("O6").Value = [double]"0.9999999999999872"
$ws.Range("P6").Value = [double]"0.9116650907698488"
$ws.Range("Q6").Value = [double]"31.26691834756943"
$ws.Range("R6").Value = [double]"-89.9999999999968"
$ws.Range("S6").Value = [double]"141.6179324093825"

# LLG_max_6 (sheet25.xml)
$ws = $wb.Worksheets.Item("LLG_max_6")
$ws.Range("N2").Value = [double]"1.039340378446755"
$ws.Range("O2").Value = [double]"0.9380319009586665"
$ws.Range("P2").Value = [double]"1.040927986637023"
$ws.Range("Q2").Value = [double]"25.58585507467569"
$ws.Range("R2").Value = [double]"-91.13019952199208"
$ws.Range("S2").Value = [double]"151.9806373353579"
$ws.Range("N3").Value = [double]"0.9011649575206414"
$ws.Range("O3").Value = [double]"0.291705603751404"
$ws.Range("P3").Value = [double]"0.795729159302157"
$ws.Range("Q3").Value = [double]"0.9710074319138136"
$ws.Range("R3").Value = [double]"-119.1713773546608"
$ws.Range("S3").Value = [double]"162.4882680064864"
$ws.Range("C4").Value = [double]"4.649268190431195"
$ws.Range("D4").Value = [double]"3.168041376878388"
$ws.Range("F4").Value = [double]"53.68512482560429"
$ws.Range("G4").Value = [double]"36.5813908348922"
$ws.Range("H4").Value = [double]"4.37704938291617"
$ws.Range("I4").Value = [double]"3.22760593527366"
$ws.Range("J4").Value = [double]"0.6277319163569961"
$ws.Range("K4").Value = [double]"2.898438697948067"
$ws.Range("L4").Value = [double]"0.627731916364051"
$ws.Range("M4").Value = [double]"2.898438697960011"
$ws.Range("N4").Value = [double]"0.7821187621557867"
$ws.Range("P4").Value = [double]"0.7821187621601789"
$ws.Range("Q4").Value = [double]"-8.509014703356241"
$ws.Range("S4").Value = [double]"171.4909852966517"
$ws.Range("T4").Value = [double]"2.876291984762022"
$ws.Range("N5").Value = [double]"0.7821187621457937"
$ws.Range("P5").Value = [double]"0.7821187621701671"
$ws.Range("Q5").Value = [double]"-8.509014702545327"
$ws.Range("S5").Value = [double]"171.4909852965164"
$ws.Range("N6").Value = [double]"0.7821187621457937"
$ws.Range("P6").Value = [double]"0.7821187621701671"
$ws.Range("Q6").Value = [double]"-8.509014702545327"
$ws.Range("S6").Value = [double]"171.4909852965164"

# LLG_max_10 (sheet26.xml)
$ws = $wb.Worksheets.Item("LLG_max_10")
$ws.Range("N2").Value = [double]"1.039340378446755"
$ws.Range("O2").Value = [double]"0.9380319009586665"
$ws.Range("P2").Value = [double]"1.040927986637023"
$ws.Range("Q2").Value = [double]"25.58585507467569"
$ws.Range("R2").Value = [double]"-91.13019952199208"
$ws.Range("S2").Value = [double]"151.9806373353579"
$ws.Range("N3").Value = [double]"0.9011649575206414"
$ws.Range("O3").Value = [double]"0.291705603751404"
$ws.Range("P3").Value = [double]"0.795729159302157"
$ws.Range("Q3").Value = [double]"0.9710074319138136"
$ws.Range("R3").Value = [double]"-119.1713773546608"
$ws.Range("S3").Value = [double]"162.4882680064864"
$ws.Range("C4").Value = [double]"4.649268190431195"
$ws.Range("D4").Value = [double]"3.168041376878388"
$ws.Range("F4").Value = [double]"53.68512482560429"
$ws.Range("G4").Value = [double]"36.5813908348922"
$ws.Range("H4").Value = [double]"4.37704938291617"
$ws.Range("I4").Value = [double]"3.22760593527366"
$ws.Range("J4").Value = [double]"0.6277319163569961"
$ws.Range("K4").Value = [double]"2.898438697948067"
$ws.Range("L4").Value = [double]"0.627731916364051"
$ws.Range("M4").Value = [double]"2.898438697960011"
$ws.Range("N4").Value = [double]"0.7821187621557867"
$ws.Range("P4").Value = [double]"0.7821187621601789"
$ws.Range("Q4").Value = [double]"-8.509014703356241"
$ws.Range("S4").Value = [double]"171.4909852966517"
$ws.Range("T4").Value = [double]"2.876291984762022"
$ws.Range("N5").Value = [double]"0.7821187621457937"
$ws.Range("P5").Value = [double]"0.7821187621701671"
$ws.Range("Q5").Value = [double]"-8.509014702545327"
$ws.Range("S5").Value = [double]"171.4909852965164"
$ws.Range("N6").Value = [double]"0.7821187621457937"
$ws.Range("P6").Value = [double]"0.7821187621701671"
$ws.Range("Q6").Value = [double]"-8.509014702545327"
$ws.Range("S6").Value = [double]"171.4909852965164"

# LLG_max_fault_6 (sheet27.xml)
$ws = $wb.Worksheets.Item("LLG_max_fault_6")
$ws.Range("N2").Value = [double]"1.039340378446755"
$ws.Range("O2").Value = [double]"0.9380319009586665"
$ws.Range("P2").Value = [double]"1.040927986637023"
$ws.Range("Q2").Value = [double]"25.58585507467569"
$ws.Range("R2").Value = [double]"-91.13019952199208"
$ws.Range("S2").Value = [double]"151.9806373353579"
$ws.Range("N3").Value = [double]"0.9011649575206414"
$ws.Range("O3").Value = [double]"0.291705603751404"
$ws.Range("P3").Value = [double]"0.795729159302157"
$ws.Range("Q3").Value = [double]"0.9710074319138136"
$ws.Range("R3").Value = [double]"-119.1713773546608"
$ws.Range("S3").Value = [double]"162.4882680064864"
$ws.Range("C4").Value = [double]"4.649268190431195"
$ws.Range("D4").Value = [double]"3.168041376878388"
$ws.Range("F4").Value = [double]"53.68512482560429"
$ws.Range("G4").Value = [double]"36.5813908348922"
$ws.Range("H4").Value = [double]"4.37704938291617"
$ws.Range("I4").Value = [double]"3.22760593527366"
$ws.Range("J4").Value = [double]"0.6277319163569961"
$ws.Range("K4").Value = [double]"2.898438697948067"
$ws.Range("L4").Value = [double]"0.627731916364051"
$ws.Range("M4").Value = [double]"2.898438697960011"
$ws.Range("N4").Value = [double]"0.7821187621557867"
$ws.Range("P4").Value = [double]"0.7821187621601789"
$ws.Range("Q4").Value = [double]"-8.509014703356241"
$ws.Range("S4").Value = [double]"171.4909852966517"
$ws.Range("T4").Value = [double]"2.876291984762022"
$ws.Range("N5").Value = [double]"0.7821187621457937"
$ws.Range("P5").Value = [double]"0.7821187621701671"
$ws.Range("Q5").Value = [double]"-8.509014702545327"
$ws.Range("S5").Value = [double]"171.4909852965164"
$ws.Range("N6").Value = [double]"0.7821187621457937"
$ws.Range("P6").Value = [double]"0.7821187621701671"
$ws.Range("Q6").Value = [double]"-8.509014702545327"
$ws.Range("S6").Value = [double]"171.4909852965164"

# LLG_max_fault_10 (sheet28.xml)
$ws = $wb.Worksheets.Item("LLG_max_fault_10")
$ws.Range("N2").Value = [double]"1.039340378446755"
$ws.Range("O2").Value = [double]"0.9380319009586665"
$ws.Range("P2").Value = [double]"1.040927986637023"
$ws.Range("Q2").Value = [double]"25.58585507467569"
$ws.Range("R2").Value = [double]"-91.13019952199208"
$ws.Range("S2").Value = [double]"151.9806373353579"
$ws.Range("N3").Value = [double]"0.9011649575206414"
$ws.Range("O3").Value = [double]"0.291705603751404"
$ws.Range("P3").Value = [double]"0.795729159302157"
$ws.Range("Q3").Value = [double]"0.9710074319138136"
$ws.Range("R3").Value = [double]"-119.1713773546608"
$ws.Range("S3").Value = [double]"162.4882680064864"
$ws.Range("C4").Value = [double]"4.649268190431195"
$ws.Range("D4").Value = [double]"3.168041376878388"
$ws.Range("F4").Value = [double]"53.68512482560429"
$ws.Range("G4").Value = [double]"36.5813908348922"
$ws.Range("H4").Value = [double]"4.37704938291617"
$ws.Range("I4").Value = [double]"3.22760593527366"
$ws.Range("J4").Value = [double]"0.6277319163569961"
$ws.Range("K4").Value = [double]"2.898438697948067"
$ws.Range("L4").Value = [double]"0.627731916364051"
$ws.Range("M4").Value = [double]"2.898438697960011"
$ws.Range("N4").Value = [double]"0.7821187621557867"
$ws.Range("P4").Value = [double]"0.7821187621601789"
$ws.Range("Q4").Value = [double]"-8.509014703356241"
$ws.Range("S4").Value = [double]"171.4909852966517"
$ws.Range("T4").Value = [double]"2.876291984762022"
$ws.Range("N5").Value = [double]"0.7821187621457937"
$ws.Range("P5").Value = [double]"0.7821187621701671"
$ws.Range("Q5").Value = [double]"-8.509014702545327"
$ws.Range("S5").Value = [double]"171.4909852965164"
$ws.Range("N6").Value = [double]"0.7821187621457937"
$ws.Range("P6").Value = [double]"0.7821187621701671"
$ws.Range("Q6").Value = [double]"-8.509014702545327"
$ws.Range("S6").Value = [double]"171.4909852965164"

# LLG_min_6 (sheet29.xml)
$ws = $wb.Worksheets.Item("LLG_min_6")
$ws.Range("N2").Value = [double]"0.9578472943291766"
$ws.Range("O2").Value = [double]"0.8500049525399896"
$ws.Range("P2").Value = [double]"0.9405846205244769"
$ws.Range("Q2").Value = [double]"24.98082684303258"
$ws.Range("R2").Value = [double]"-92.6538163133243"
$ws.Range("S2").Value = [double]"151.7926647951937"
$ws.Range("N3").Value = [double]"0.9015751892214467"
$ws.Range("O3").Value = [double]"0.3617201835036115"
$ws.Range("P3").Value = [double]"0.7146959486074793"
$ws.Range("Q3").Value = [double]"2.491133493871041"
$ws.Range("R3").Value = [double]"-129.0915432241699"
$ws.Range("S3").Value = [double]"160.2456697546928"
$ws.Range("C4").Value = [double]"3.820497609465972"
$ws.Range("D4").Value = [double]"2.668597153164733"
$ws.Range("F4").Value = [double]"44.11530646527001"
$ws.Range("G4").Value = [double]"30.81430569476655"
$ws.Range("H4").Value = [double]"6.086106992372931"
$ws.Range("I4").Value = [double]"3.212307515692304"
$ws.Range("J4").Value = [double]"1.113751880246063"
$ws.Range("K4").Value = [double]"3.006397543767953"
$ws.Range("L4").Value = [double]"1.113751880233333"
$ws.Range("M4").Value = [double]"3.006397543745581"
$ws.Range("N4").Value = [double]"0.7315844258515088"
$ws.Range("P4").Value = [double]"0.7315844258518936"
$ws.Range("Q4").Value = [double]"-7.543414800613841"
$ws.Range("S4").Value = [double]"172.4565851989179"
$ws.Range("T4").Value = [double]"2.126132628031408"
$ws.Range("N5").Value = [double]"0.7315844258432987"
$ws.Range("P5").Value = [double]"0.7315844258633758"
$ws.Range("Q5").Value = [double]"-7.543414799225709"
$ws.Range("S5").Value = [double]"172.4565851983412"
$ws.Range("N6").Value = [double]"0.7315844258432987"
$ws.Range("P6").Value = [double]"0.7315844258633758"
$ws.Range("Q6").Value = [double]"-7.543414799225709"
$ws.Range("S6").Value = [double]"172.4565851983412"

# LLG_min_10 (sheet30.xml)
$ws = $wb.Worksheets.Item("LLG_min_10")
$ws.Range("N2").Value = [double]"0.9578472943291766"
$ws.Range("O2").Value = [double]"0.8500049525399896"
$ws.Range("P2").Value = [double]"0.9405846205244769"
$ws.Range("Q2").Value = [double]"24.98082684303258"
$ws.Range("R2").Value = [double]"-92.6538163133243"
$ws.Range("S2").Value = [double]"151.7926647951937"
$ws.Range("N3").Value = [double]"0.9015751892214467"
$ws.Range("O3").Value = [double]"0.3617201835036115"
$ws.Range("P3").Value = [double]"0.7146959486074793"
$ws.Range("Q3").Value = [double]"2.491133493871041"
$ws.Range("R3").Value = [double]"-129.0915432241699"
$ws.Range("S3").Value = [double]"160.2456697546928"
$ws.Range("C4").Value = [double]"3.820497609465972"
$ws.Range("D4").Value = [double]"2.668597153164733"
$ws.Range("F4").Value = [double]"44.11530646527001"
$ws.Range("G4").Value = [double]"30.81430569476655"
$ws.Range("H4").Value = [double]"6.086106992372931"
$ws.Range("I4").Value = [double]"3.212307515692304"
$ws.Range("J4").Value = [double]"1.113751880246063"
$ws.Range("K4").Value = [double]"3.006397543767953"
$ws.Range("L4").Value = [double]"1.113751880233333"
$ws.Range("M4").Value = [double]"3.006397543745581"
$ws.Range("N4").Value = [double]"0.7315844258515088"
$ws.Range("P4").Value = [double]"0.7315844258518936"
$ws.Range("Q4").Value = [double]"-7.543414800613841"
$ws.Range("S4").Value = [double]"172.4565851989179"
$ws.Range("T4").Value = [double]"2.126132628031408"
$ws.Range("N5").Value = [double]"0.7315844258432987"
$ws.Range("P5").Value = [double]"0.7315844258633758"
$ws.Range("Q5").Value = [double]"-7.543414799225709"
$ws.Range("S5").Value = [double]"172.4565851983412"
$ws.Range("N6").Value = [double]"0.7315844258432987"
$ws.Range("P6").Value = [double]"0.7315844258633758"
$ws.Range("Q6").Value = [double]"-7.543414799225709"
$ws.Range("S6").Value = [double]"172.4565851983412"

# LLG_min_fault_6 (sheet31.xml)
$ws = $wb.Worksheets.Item("LLG_min_fault_6")
$ws.Range("N2").Value = [double]"0.9578472943291766"
$ws.Range("O2").Value = [double]"0.8500049525399896"
$ws.Range("P2").Value = [double]"0.9405846205244769"
$ws.Range("Q2").Value = [double]"24.98082684303258"
$ws.Range("R2").Value = [double]"-92.6538163133243"
$ws.Range("S2").Value = [double]"151.7926647951937"
$ws.Range("N3").Value = [double]"0.9015751892214467"
$ws.Range("O3").Value = [double]"0.3617201835036115"
$ws.Range("P3").Value = [double]"0.7146959486074793"
$ws.Range("Q3").Value = [double]"2.491133493871041"
$ws.Range("R3").Value = [double]"-129.0915432241699"
$ws.Range("S3").Value = [double]"160.2456697546928"
$ws.Range("C4").Value = [double]"3.820497609465972"
$ws.Range("D4").Value = [double]"2.668597153164733"
$ws.Range("F4").Value = [double]"44.11530646527001"
$ws.Range("G4").Value = [double]"30.81430569476655"
$ws.Range("H4").Value = [double]"6.086106992372931"
$ws.Range("I4").Value = [double]"3.212307515692304"
$ws.Range("J4").Value = [double]"1.113751880246063"
$ws.Range("K4").Value = [double]"3.006397543767953"
$ws.Range("L4").Value = [double]"1.113751880233333"
$ws.Range("M4").Value = [double]"3.006397543745581"
$ws.Range("N4").Value = [double]"0.7315844258515088"
$ws.Range("P4").Value = [double]"0.7315844258518936"
$ws.Range("Q4").Value = [double]"-7.543414800613841"
$ws.Range("S4").Value = [double]"172.4565851989179"
$ws.Range("T4").Value = [double]"2.126132628031408"
$ws.Range("N5").Value = [double]"0.7315844258432987"
$ws.Range("P5").Value = [double]"0.7315844258633758"
$ws.Range("Q5").Value = [double]"-7.543414799225709"
$ws.Range("S5").Value = [double]"172.4565851983412"
$ws.Range("N6").Value = [double]"0.7315844258432987"
$ws.Range("P6").Value = [double]"0.7315844258633758"
$ws.Range("Q6").Value = [double]"-7.543414799225709"
$ws.Range("S6").Value = [double]"172.4565851983412"

# LLG_min_fault_10 (sheet32.xml)
$ws = $wb.Worksheets.Item("LLG_min_fault_10")
$ws.Range("N2").Value = [double]"0.9578472943291766"
$ws.Range("O2").Value = [double]"0.8500049525399896"
$ws.Range("P2").Value = [double]"0.9405846205244769"
$ws.Range("Q2").Value = [double]"24.98082684303258"
$ws.Range("R2").Value = [double]"-92.6538163133243"
$ws.Range("S2").Value = [double]"151.7926647951937"
$ws.Range("N3").Value = [double]"0.9015751892214467"
$ws.Range("O3").Value = [double]"0.3617201835036115"
$ws.Range("P3").Value = [double]"0.7146959486074793"
$ws.Range("Q3").Value = [double]"2.491133493871041"
$ws.Range("R3").Value = [double]"-129.0915432241699"
$ws.Range("S3").Value = [double]"160.2456697546928"
$ws.Range("C4").Value = [double]"3.820497609465972"
$ws.Range("D4").Value = [double]"2.668597153164733"
$ws.Range("F4").Value = [double]"44.11530646527001"
$ws.Range("G4").Value = [double]"30.81430569476655"
$ws.Range("H4").Value = [double]"6.086106992372931"
$ws.Range("I4").Value = [double]"3.212307515692304"
$ws.Range("J4").Value = [double]"1.113751880246063"
$ws.Range("K4").Value = [double]"3.006397543767953"
$ws.Range("L4").Value = [double]"1.113751880233333"
$ws.Range("M4").Value = [double]"3.006397543745581"
$ws.Range("N4").Value = [double]"0.7315844258515088"
$ws.Range("P4").Value = [double]"0.7315844258518936"
$ws.Range("Q4").Value = [double]"-7.543414800613841"
$ws.Range("S4").Value = [double]"172.4565851989179"
$ws.Range("T4").Value = [double]"2.126132628031408"
$ws.Range("N5").Value = [double]"0.7315844258432987"
$ws.Range("P5").Value = [double]"0.7315844258633758"
$ws.Range("Q5").Value = [double]"-7.543414799225709"
$ws.Range("S5").Value = [double]"172.4565851983412"
$ws.Range("N6").Value = [double]"0.7315844258432987"
$ws.Range("P6").Value = [double]"0.7315844258633758"
$ws.Range("Q6").Value = [double]"-7.543414799225709"
$ws.Range("S6").Value = [double]"172.4565851983412"
